$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2025-02-27 abs_activity score was recalculated -> update C106/F106
$ws.Range("C106").Value = 9.778469025197799
$ws.Range("F106").Value = 9.778469025197799

# Append the new daily scores for 2025-02-28 (rows 110-113).
# The leading apostrophe forces these date-looking strings to be stored
# as text, matching the existing Date column entries (e.g. A106 = "2025-02-27").
$ws.Range("A110").Value = "'2025-02-28"
$ws.Range("B110").Value = "abs_activity"
$ws.Range("C110").Value = 10
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 0
$ws.Range("F110").Value = 10

$ws.Range("A111").Value = "'2025-02-28"
$ws.Range("B111").Value = "rel_activity"
$ws.Range("C111").Value = 6.085069444444445
$ws.Range("D111").Value = 0
$ws.Range("E111").Value = 0
$ws.Range("F111").Value = 6.085069444444445

$ws.Range("A112").Value = "'2025-02-28"
$ws.Range("B112").Value = "abs_sleep"
$ws.Range("C112").Value = 10
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 0
$ws.Range("F112").Value = 10

$ws.Range("A113").Value = "'2025-02-28"
$ws.Range("B113").Value = "rel_sleep"
$ws.Range("C113").Value = 9.297281550165639
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 0
$ws.Range("F113").Value = 9.297281550165639
